$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) Reposition / resize the "Picture 2" (user-stories screenshot, Id=15).
#    Target EMU (from the authoritative diff):
#       off  x=173851  y=22168902
#       ext cx=10537051 cy=7866360
#    NOTE: the Shape.Left/Top/Width/Height setters round-trip the point value
#    through a single-precision float before converting to EMU, so a literal
#    "emu/12700" value can land 1 EMU short after truncation. The literals
#    below are the smallest points values whose float32 representation still
#    converts to the exact target EMU.
# ---------------------------------------------------------------------------
$pic = $s.Shapes.Item(24)
$pic.Left   = 13.689055442810059
$pic.Top    = 1745.5828857421875
$pic.Width  = 829.6890869140625
$pic.Height = 619.3984375

# ---------------------------------------------------------------------------
# 2) Add the new caption textbox under/around the picture (Id=9, "TextBox 8").
#    Passing the exact position/size straight into AddTextbox keeps full
#    double precision, so plain emu/12700 literals round-trip exactly here.
# ---------------------------------------------------------------------------
$tb = $s.Shapes.AddTextbox(1, 39.378188976377956, 1447.5, 803.9999212598425, 298.08283464566927)
$tb.Name = "TextBox 8"
$tb.Fill.Visible = $false

$tf = $tb.TextFrame
$tf.WordWrap = $true

$tr = $tf.TextRange
$tr.Text = "In our e-Voting System we have accumulated roughly  40 different use r stories so far. All of these user stories revolves around 3 actors; Voter, Admin, Mobile Voter.  Below is a snippet of just a couple of the user stories we came up with."
$tr.Font.Size = 40

# Auto-fit the box height to the text now that the text/size are in place.
$tf.AutoSize = 1
